# -----------------------------------------------------------------------------
# Refresh the cryptocurrency table on Sheet1 with the latest scraped values.
# Column D ("Price") cells that look numeric (e.g. "191.30", "7.220") are
# written with a leading apostrophe so Excel keeps them as literal text instead
# of collapsing them into a Double and dropping significant trailing/format digits.
# -----------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '26.098.78'
$ws.Range('E2').Value = '  +0.13%  '

# Row 3
$ws.Range('D3').Value = '1.653.06'
$ws.Range('E3').Value = '  +0.21%  '

# Row 4
$ws.Range('D4').Value = '''1.003'
$ws.Range('E4').Value = '  -0.38%  '

# Row 5
$ws.Range('D5').Value = '''218.44'
$ws.Range('E5').Value = '  +0.46%  '

# Row 6
$ws.Range('D6').Value = '''0.5202'
$ws.Range('E6').Value = '  -0.13%  '

# Row 7
$ws.Range('E7').Value = '  -0.39%  '

# Row 8
$ws.Range('D8').Value = '''0.2646'
$ws.Range('E8').Value = '  +1.22%  '

# Row 9
$ws.Range('D9').Value = '''0.06336'
$ws.Range('E9').Value = '  +0.86%  '

# Row 10
$ws.Range('D10').Value = '''20.41'
$ws.Range('E10').Value = '  -0.31%  '

# Row 11
$ws.Range('D11').Value = '''0.07692'

# Row 12
$ws.Range('D12').Value = '''4.605'
$ws.Range('E12').Value = '  +2.88%  '

# Row 13
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '1.880.18'
$ws.Range('E13').Value = '  +0.14%  '

# Row 14
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.602.35'
$ws.Range('E14').Value = '  -2.78%  '

# Row 15
$ws.Range('D15').Value = '''0.5593'
$ws.Range('E15').Value = '  +1.12%  '

# Row 16
$ws.Range('D16').Value = '0.0₅8161'
$ws.Range('E16').Value = '  +2.17%  '

# Row 17
$ws.Range('D17').Value = '''65.38'
$ws.Range('E17').Value = '  +1.00%  '

# Row 18
$ws.Range('D18').Value = '26.115.45'
$ws.Range('E18').Value = '  +0.21%  '

# Row 19
$ws.Range('E19').Value = '  -0.33%  '

# Row 20
$ws.Range('D20').Value = '''4.631'
$ws.Range('E20').Value = '  +0.24%  '

# Row 21
$ws.Range('D21').Value = '''10.49'
$ws.Range('E21').Value = '  +4.38%  '

# Row 22
$ws.Range('D22').Value = '''191.30'
$ws.Range('E22').Value = '  -1.35%  '

# Row 23
$ws.Range('D23').Value = '''5.925'
$ws.Range('E23').Value = '  -0.26%  '

# Row 24
$ws.Range('D24').Value = '''1.004'
$ws.Range('E24').Value = '  -0.27%  '

# Row 25
$ws.Range('D25').Value = '''144.17'
$ws.Range('E25').Value = '  -1.69%  '

# Row 26
$ws.Range('D26').Value = '''0.1189'
$ws.Range('E26').Value = '  -0.96%  '

# Row 27
$ws.Range('D27').Value = '''7.220'
$ws.Range('E27').Value = '  +0.68%  '

# Row 28
$ws.Range('D28').Value = '''15.90'
$ws.Range('E28').Value = '  +0.17%  '

# Row 29
$ws.Range('D29').Value = '''1.503'
$ws.Range('E29').Value = '  +1.68%  '

# Row 30
$ws.Range('D30').Value = '''0.05485'
$ws.Range('E30').Value = '  -1.99%  '

# Row 31
$ws.Range('D31').Value = '''1.269'
$ws.Range('E31').Value = '  +0.32%  '

# Row 32
$ws.Range('D32').Value = '''3.445'
$ws.Range('E32').Value = '  -0.99%  '

# Row 33
$ws.Range('D33').Value = '''3.354'
$ws.Range('E33').Value = '  -0.06%  '

# Row 34
$ws.Range('D34').Value = '''1.560'
$ws.Range('E34').Value = '  -1.72%  '

# Row 35
$ws.Range('E35').Value = '  +0.38%  '

# Row 36
$ws.Range('D36').Value = '''0.9475'
$ws.Range('E36').Value = '  +0.03%  '

# Row 37
$ws.Range('D37').Value = '''2.785'
$ws.Range('E37').Value = '  -0.41%  '

# Row 38
$ws.Range('D38').Value = '''0.5643'
$ws.Range('E38').Value = '  +0.25%  '

# Row 39
$ws.Range('D39').Value = '''0.01580'
$ws.Range('E39').Value = '  -0.18%  '

# Row 40
$ws.Range('D40').Value = '''5.847'
$ws.Range('E40').Value = '  -1.65%  '

# Row 42
$ws.Range('D42').Value = '1.027.11'
$ws.Range('E42').Value = '  -2.89%  '

# Row 43
$ws.Range('D43').Value = '''0.8284'

# Row 44
$ws.Range('D44').Value = '''101.21'
$ws.Range('E44').Value = '  -1.06%  '

# Row 45
$ws.Range('D45').Value = '1.795.15'
$ws.Range('E45').Value = '  +0.35%  '

# Row 46
$ws.Range('D46').Value = '''57.59'
$ws.Range('E46').Value = '  +0.89%  '

# Row 47
$ws.Range('E47').Value = '  +7.11%  '

# Row 48
$ws.Range('D48').Value = '''0.9986'
$ws.Range('E48').Value = '  -0.86%  '

# Row 49
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').Value = '''0.4336'
$ws.Range('E49').Value = '  +0.13%  '

# Row 50
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '''7.973'
$ws.Range('E50').Value = '  +0.44%  '

# Row 51
$ws.Range('D51').Value = '''0.05168'
$ws.Range('E51').Value = '  -2.67%  '
